$wb = $excel.ActiveWorkbook

function Set-ObjTablesHeader {
    param($ws, $cellRef, $newValue)
    $ws.Unprotect()
    $ws.Range($cellRef).Value = $newValue
    $ws.Protect()
}

$ws = $wb.Worksheets.Item(1)
Set-ObjTablesHeader $ws "A1" "!!!ObjTables objTablesVersion='0.0.8' date='2020-03-09 13:01:36'"
$ws = $wb.Worksheets.Item(1)
Set-ObjTablesHeader $ws "A2" "!!ObjTables type='Data' id='Compartment' name='Compartment' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(2)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Compound' name='Compound' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(3)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Definition' name='Definition' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(4)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Enzyme' name='Enzyme' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(5)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='FbcObjective' name='FbcObjective' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(6)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Gene' name='Gene' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(7)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Layout' name='Layout' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(8)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Measurement' name='Measurement' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(9)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='PbConfig' name='PbConfig' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(10)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Position' name='Position' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(11)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Protein' name='Protein' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(12)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Quantity' name='Quantity' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(13)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='QuantityInfo' name='QuantityInfo' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(14)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='QuantityMatrix' name='QuantityMatrix' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(15)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Reaction' name='Reaction' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(16)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='ReactionStoichiometry' name='ReactionStoichiometry' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(17)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Regulator' name='Regulator' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(18)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Relation' name='Relation' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(19)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='Relationship' name='Relationship' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(20)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='SparseMatrix' name='SparseMatrix' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(21)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='SparseMatrixColumn' name='SparseMatrixColumn' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(22)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='SparseMatrixOrdered' name='SparseMatrixOrdered' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(23)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='SparseMatrixRow' name='SparseMatrixRow' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(24)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='StoichiometricMatrix' name='StoichiometricMatrix' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(25)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='rxnconContingencyList' name='rxnconContingencyList' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"
$ws = $wb.Worksheets.Item(26)
Set-ObjTablesHeader $ws "A1" "!!ObjTables type='Data' id='rxnconReactionList' name='rxnconReactionList' date='2020-03-09 13:01:36' objTablesVersion='0.0.8'"

Write-Output "done"